# Regenerate the "K" column (column G, header "K" -- strikeouts) of the
# per-game save_data sheet. The sheet previously stored a different metric
# ("Strike#") in column G; this replaces each value with the real K value,
# matching a freshly regenerated pull of the underlying stats (std/mean and
# s_vals are derived downstream from this data, outside this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for data rows 2..65 (one entry per row, in row order).
$newK = @(
    0, 2, 0, 1, 1, 1, 1, 0, 0, 1,
    1, 1, 1, 2, 0, 3, 1, 1, 2, 1,
    2, 1, 2, 2, 1, 0, 1, 1, 2, 1,
    0, 0, 1, 1, 1, 1, 1, 1, 1, 1,
    1, 2, 1, 1, 1, 2, 1, 1, 1, 1,
    0, 1, 1, 1, 1, 0, 1, 1, 0, 4,
    0, 2, 2, 1
)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
